# The deck currently uses the "Integral" theme (ppt/theme/theme2.xml,
# wired to the slide master / presentation) while ppt/theme/theme1.xml
# (wired only to the notes master) carries the default "Office Theme"
# palette. The authored change swaps the two: the presentation's visible
# design becomes the plain "Office Theme" colour palette.
#
# Helper: convert an RRGGBB hex string into the OLE-packed 0x00BBGGRR
# integer PowerPoint's COM object model expects for ColorFormat.RGB /
# ThemeColor.RGB.
function HexToRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# The presentation's single Design -> SlideMaster -> Theme carries the
# colour scheme that is persisted back into the theme XML part on save.
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Target palette: the standard Office theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) in ThemeColorScheme index order.
$officeThemeHex = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToRGB $officeThemeHex[$i - 1]
}
